$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must be forced to text to preserve exact
# formatting (e.g. trailing zeros, multi-dot thousand separators) --
# otherwise Excel auto-converts numeric-looking strings to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "74.196.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.625.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "186.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "582.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.197"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.634.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "74.318.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.119.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000188"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.646.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0935"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "524.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "162.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0847"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.521"
$ws.Range("D51").Style = "Normal"

# Coin name / link / volume columns (already safe as text)
$ws.Range("E2").Value = "  +6.06%  "
$ws.Range("E3").Value = "  +6.37%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("E5").Value = "  +11.94%  "
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.43%  "
$ws.Range("E9").Value = "  +11.58%  "
$ws.Range("E10").Value = "  +6.76%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +5.99%  "
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("E14").Value = "  +6.41%  "
$ws.Range("E15").Value = "  +6.79%  "
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("E17").Value = "  +9.95%  "
$ws.Range("E18").Value = "  +7.27%  "
$ws.Range("E19").Value = "  +27.84%  "
$ws.Range("E20").Value = "  +8.63%  "
$ws.Range("E21").Value = "  +6.18%  "
$ws.Range("E22").Value = "  +12.78%  "
$ws.Range("E23").Value = "  +4.72%  "
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +5.45%  "
$ws.Range("E27").Value = "  +4.78%  "
$ws.Range("E28").Value = "  +8.41%  "
$ws.Range("E29").Value = "  +6.51%  "
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("E31").Value = "  +9.49%  "
$ws.Range("E32").Value = "  +16.36%  "
$ws.Range("E33").Value = "  +9.77%  "
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("E35").Value = "  +5.96%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("E38").Value = "  +7.04%  "
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +8.68%  "
$ws.Range("E43").Value = "  +7.21%  "
$ws.Range("E44").Value = "  +6.24%  "
$ws.Range("E45").Value = "  +22.09%  "
$ws.Range("E46").Value = "  +9.59%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E48").Value = "  +6.49%  "
$ws.Range("E49").Value = "  +16.83%  "
$ws.Range("E50").Value = "  +5.40%  "
$ws.Range("E51").Value = "  +5.82%  "
